$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the formatting
# (bold/border/centered header style) already used by H1 and the rest of
# row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"

# Fill the new I (I0) and J (IF) columns for every data row: I is always 1,
# J duplicates the existing H (IP) value for that row.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 9).Value2 = 1
    $hValue = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 10).Value2 = $hValue
}
